$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.258.80"
$ws.Range("E2").Value = "  +2.19%  "
$ws.Range("D3").Value = "1.996.42"
$ws.Range("E3").Value = "  +6.04%  "
$ws.Range("D4").Value = "'0.9990"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'0.7792"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +64.49%  "
$ws.Range("D6").Value = "'254.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.47%  "
$ws.Range("D7").Value = "'0.9990"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'0.3474"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +20.30%  "
$ws.Range("D9").Value = "'28.03"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +25.67%  "
$ws.Range("D10").Value = "'0.07085"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.52%  "
$ws.Range("D11").Value = "'0.8424"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +11.27%  "
$ws.Range("D12").Value = "'0.08192"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.67%  "
$ws.Range("D13").Value = "'100.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("D14").Value = "1.995.80"
$ws.Range("E14").Value = "  +6.01%  "
$ws.Range("D15").Value = "'5.641"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.79%  "
$ws.Range("D16").Value = "'15.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +16.51%  "
$ws.Range("D17").Value = "'272.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.75%  "
$ws.Range("D18").Value = "31.245.18"
$ws.Range("D19").Value = "'5.987"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +12.10%  "
$ws.Range("D20").Value = "'0.000008026"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.76%  "
$ws.Range("D21").Value = "2.257.96"
$ws.Range("E21").Value = "  +6.26%  "
$ws.Range("D22").Value = "'0.9984"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Value = "'0.9987"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "'7.107"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +10.76%  "
$ws.Range("D25").Value = "'10.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.33%  "
$ws.Range("D26").Value = "'164.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("D27").Value = "'0.1432"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +46.91%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'19.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.01%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'2.410"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +26.50%  "
$ws.Range("D30").Value = "'1.593"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.95%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.635"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.32%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'1.365"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.04%  "
$ws.Range("D33").Value = "'4.455"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.59%  "
$ws.Range("D34").Value = "'0.05341"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.31%  "
$ws.Range("D35").Value = "'1.257"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.06%  "
$ws.Range("D36").Value = "'0.7910"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +13.31%  "
$ws.Range("D37").Value = "'2.767"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "'0.9985"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("D39").Value = "'0.02007"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.49%  "
$ws.Range("D40").Value = "'2.919"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.74%  "
$ws.Range("D41").Value = "'84.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.75%  "
$ws.Range("D42").Value = "'6.791"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.86%  "
$ws.Range("D43").Value = "'0.4689"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.73%  "
$ws.Range("D44").Value = "'2.140"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.83%  "
$ws.Range("D45").Value = "'0.8579"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.12%  "
$ws.Range("D46").Value = "'105.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.99%  "
$ws.Range("D47").Value = "'0.9992"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("D48").Value = "'7.784"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.06%  "
$ws.Range("D49").Value = "'10.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").Value = "'37.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.42%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "'3.009"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +45.75%  "
